$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the latest daily profit entry (run on 2025-10-04) as row 48,
# continuing the existing Date/Profit table that currently ends at row 47.
$row = 48

# Force text formatting first so the date-like string isn't auto-converted
# into a date serial number (matches the plain text style of prior rows),
# then reset the cell style back to Normal so no extra style index lingers.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "10/04/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 15102.67
